# Edit: in the "Segmentación clara" bullet, drop "-bollería" from
# "pastelería-bollería" and insert a new clause " en un menú desplegable
# en la sección productos" right before ", reflejando la organización
# del negocio." — splitting the original single run into three runs
# (same rPr) the way Word does when text is typed/edited in separate
# passes.

$d = $word.ActiveDocument

# 1) "pastelería-bollería" -> "pastelería"
$d.Content.Find.Execute(
    "pastelería-bollería y empanadas)", $true, $false, $false, $false,
    $false, $true, 1, $false, "pastelería y empanadas)", 2
)

# 2) Insert the new clause right before the trailing ", reflejando..." text
$tail = $d.Content
$tail.Find.Execute(", reflejando la organización del negocio.")
$tail.InsertBefore(" en un menú desplegable en la sección productos")

# 3) Force the boundary between "...empanadas)" and " en un menú..." to
#    become a distinct run (Word keeps these as separate <w:r> elements
#    even though the formatting is identical) by toggling a character
#    property on/off over the inserted+trailing text.
$afterParen = $d.Content
$afterParen.Find.Execute(" en un menú desplegable en la sección productos, reflejando la organización del negocio.")
$afterParen.Bold = 1
$afterParen.Bold = 0

# 4) Likewise force the boundary between "...productos" and ", reflejando..."
$beforeComma = $d.Content
$beforeComma.Find.Execute(", reflejando la organización del negocio.")
$beforeComma.Bold = 1
$beforeComma.Bold = 0
